$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column F ("dSF") per row, as computed/repulled data.
$updates = @{
    3  = -1
    4  = -1
    5  = 3
    6  = 1
    7  = -3
    8  = 8
    10 = 0
    11 = 9
    12 = -3
    13 = -3
    14 = 4
    15 = -1
    17 = 1
    18 = 7
    19 = 0
    20 = -3
    21 = 2
    22 = -3
    23 = 4
    24 = 7
    25 = 1
    26 = -1
    27 = 2
    28 = -1
    30 = 2
    32 = -2
    33 = 2
    34 = 2
    35 = -2
    36 = -3
    37 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
